$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows in the Tableau section (after TableauPile7 at row 25,
# before the blank separator that is currently row 26) for the new
# "IsOpen"/"IsClosed" bool properties. This pushes the existing blank
# separator row (which carries row-level custom formatting) down from 26 to
# 28, and everything below (Foundation section, Board section) down by two
# rows as well.
$ws.Rows("26:27").Insert()

# Each TableauPile1..7 row (19-25) gets a "List<Card>" type in column B.
$ws.Range("B19").Value = "List<Card>"
$ws.Range("B20").Value = "List<Card>"
$ws.Range("B21").Value = "List<Card>"
$ws.Range("B22").Value = "List<Card>"
$ws.Range("B23").Value = "List<Card>"
$ws.Range("B24").Value = "List<Card>"
$ws.Range("B25").Value = "List<Card>"

# New row 26: IsOpen / bool
$ws.Range("A26").Value = "IsOpen"
$ws.Range("B26").Value = "bool"

# New row 27: IsClosed / bool -- matches the formatting used elsewhere in
# this section (same number format as the blank separator row / the
# FoundationPile rows below, style index carrying numFmtId 15).
$ws.Range("A27").Value = "IsClosed"
$ws.Range("B27").Value = "bool"
$ws.Range("A27:B27").NumberFormat = "d-mmm-yy"

# Update the view: scroll position and active selection moved.
$ws.Range("P15").Select()
$excel.ActiveWindow.ScrollRow = 6
